# Applies the "advanced feasibility notes with more specific guidelines" edit:
#  - Results sheet: the optimizer now finds a feasible allocation (previously
#    every Bid ID showed "No Bid" / all-zero economics). Rows 2-10 get their
#    award columns (F,G,H,J,K,L,M) populated, a couple of Bid ID/Facility/
#    Baseline values shift because Bid ID 3 now splits across two rows, and a
#    new split row is inserted at row 11 (old row 11 shifts down to row 12).
#  - Feasibility Notes sheet: the rule-evaluation note text is replaced.
#  - LP Model sheet: the bid-exclusion constraints are removed and a
#    '# of Transitions' rule constraint is added.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Results sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Results")

# Insert a new row at 11 -- this pushes the existing row 11 (Bid ID 10 /
# Facility 10 / baseline 231) down to row 12, matching the diff where a new
# row ("Bid ID Split" B for Bid ID 3, and a new split for Bid ID 10 /
# baseline 234) is inserted ahead of it.
$ws.Rows.Item(11).Insert()

# Row 2 (Bid ID 1) -- awarded to Supplier A
$ws.Cells.Item(2,6).Value  = 70000
$ws.Cells.Item(2,7).Value  = "A"
$ws.Cells.Item(2,8).Value  = 0.5252
$ws.Cells.Item(2,10).Value = 0.5252
$ws.Cells.Item(2,11).Value = 367.64
$ws.Cells.Item(2,12).Value = 700
$ws.Cells.Item(2,13).Value = 69632.36

# Row 3 (Bid ID 2) -- awarded supplier changes from A to B
$ws.Cells.Item(3,7).Value = "B"

# Row 4 (Bid ID 3, first split row) -- awarded to Supplier B
$ws.Cells.Item(4,6).Value  = 219537
$ws.Cells.Item(4,7).Value  = "B"
$ws.Cells.Item(4,8).Value  = 65
$ws.Cells.Item(4,10).Value = 65
$ws.Cells.Item(4,11).Value = 33735
$ws.Cells.Item(4,12).Value = 519
$ws.Cells.Item(4,13).Value = 185802

# Row 5 -- becomes the second split row of Bid ID 3 (Bid ID Split "B"),
# with Facility/baseline shifted from what used to be Bid ID 4's row
$ws.Cells.Item(5,1).Value  = 3
$ws.Cells.Item(5,2).Value  = "B"
$ws.Cells.Item(5,5).Value  = 423
$ws.Cells.Item(5,6).Value  = 34263
$ws.Cells.Item(5,7).Value  = "C"
$ws.Cells.Item(5,8).Value  = 60
$ws.Cells.Item(5,10).Value = 60
$ws.Cells.Item(5,11).Value = 4860
$ws.Cells.Item(5,12).Value = 81
$ws.Cells.Item(5,13).Value = 29403

# Row 6 -- now Bid ID 4 / Facility 4 (shifted up from the old Bid ID 5 /
# Facility 5 row), awarded to Supplier A
$ws.Cells.Item(6,1).Value  = 4
$ws.Cells.Item(6,3).Value  = "Facility 4"
$ws.Cells.Item(6,5).Value  = 453
$ws.Cells.Item(6,6).Value  = 2568510
$ws.Cells.Item(6,7).Value  = "A"
$ws.Cells.Item(6,8).Value  = 23
$ws.Cells.Item(6,10).Value = 23
$ws.Cells.Item(6,11).Value = 130410
$ws.Cells.Item(6,12).Value = 5670
$ws.Cells.Item(6,13).Value = 2438100

# Row 7 -- now Bid ID 5 / Facility 5 (shifted up from the old Bid ID 6 /
# Facility 9 row), awarded to Supplier B
$ws.Cells.Item(7,1).Value  = 5
$ws.Cells.Item(7,3).Value  = "Facility 5"
$ws.Cells.Item(7,5).Value  = 342
$ws.Cells.Item(7,6).Value  = 15390
$ws.Cells.Item(7,7).Value  = "B"
$ws.Cells.Item(7,8).Value  = 34
$ws.Cells.Item(7,10).Value = 34
$ws.Cells.Item(7,11).Value = 1530
$ws.Cells.Item(7,12).Value = 45
$ws.Cells.Item(7,13).Value = 13860

# Row 8 -- now Bid ID 6 (baseline shifted up from the old Bid ID 7 row),
# awarded to Supplier B
$ws.Cells.Item(8,1).Value  = 6
$ws.Cells.Item(8,5).Value  = 653
$ws.Cells.Item(8,6).Value  = 158026
$ws.Cells.Item(8,7).Value  = "B"
$ws.Cells.Item(8,8).Value  = 24
$ws.Cells.Item(8,10).Value = 24
$ws.Cells.Item(8,11).Value = 5808
$ws.Cells.Item(8,12).Value = 242
$ws.Cells.Item(8,13).Value = 152218

# Row 9 -- now Bid ID 7 (baseline shifted up from the old Bid ID 8 row),
# awarded to Supplier A
$ws.Cells.Item(9,1).Value  = 7
$ws.Cells.Item(9,5).Value  = 432
$ws.Cells.Item(9,6).Value  = 286848
$ws.Cells.Item(9,7).Value  = "A"
$ws.Cells.Item(9,8).Value  = 23
$ws.Cells.Item(9,10).Value = 23
$ws.Cells.Item(9,11).Value = 15272
$ws.Cells.Item(9,12).Value = 664
$ws.Cells.Item(9,13).Value = 271576

# Row 10 -- now Bid ID 8 / Facility 9 (shifted up from the old Bid ID 9 /
# Facility 10 row), awarded to Supplier C
$ws.Cells.Item(10,1).Value  = 8
$ws.Cells.Item(10,3).Value  = "Facility 9"
$ws.Cells.Item(10,5).Value  = 456
$ws.Cells.Item(10,6).Value  = 10944
$ws.Cells.Item(10,7).Value  = "C"
$ws.Cells.Item(10,8).Value  = 24
$ws.Cells.Item(10,10).Value = 24
$ws.Cells.Item(10,11).Value = 576
$ws.Cells.Item(10,12).Value = 24
$ws.Cells.Item(10,13).Value = 10368

# Row 11 (newly inserted) -- Bid ID 9 / Facility 10, awarded to Supplier C
# (I11/N11 hold the literal text "0%", not a formatted percentage number, so
# force Text number format first -- otherwise Excel reinterprets "0%" as the
# numeric value 0 with a Percentage number format applied.)
$ws.Cells.Item(11,9).NumberFormat  = "@"
$ws.Cells.Item(11,14).NumberFormat = "@"

$ws.Cells.Item(11,1).Value  = 9
$ws.Cells.Item(11,2).Value  = "A"
$ws.Cells.Item(11,3).Value  = "Facility 10"
$ws.Cells.Item(11,4).Value  = "C"
$ws.Cells.Item(11,5).Value  = 234
$ws.Cells.Item(11,6).Value  = 54288
$ws.Cells.Item(11,7).Value  = "C"
$ws.Cells.Item(11,8).Value  = 32
$ws.Cells.Item(11,9).Value  = "0%"
$ws.Cells.Item(11,10).Value = 32
$ws.Cells.Item(11,11).Value = 7424
$ws.Cells.Item(11,12).Value = 232
$ws.Cells.Item(11,13).Value = 46864
$ws.Cells.Item(11,14).Value = "0%"
$ws.Cells.Item(11,15).Value = 0

# Row 12 (previously row 11, pushed down by the insert) -- Bid ID 10,
# baseline 231, awarded to Supplier C
$ws.Cells.Item(12,1).Value  = 10
$ws.Cells.Item(12,6).Value  = 3003
$ws.Cells.Item(12,7).Value  = "C"
$ws.Cells.Item(12,8).Value  = 15
$ws.Cells.Item(12,10).Value = 15
$ws.Cells.Item(12,11).Value = 195
$ws.Cells.Item(12,12).Value = 13
$ws.Cells.Item(12,13).Value = 2808

# ---------------------------------------------------------------------
# 2) Feasibility Notes sheet
# ---------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Feasibility Notes")
$notesText = "Model is infeasible. Likely causes include:`n - Insufficient supplier capacity relative to demand.`n - Custom rule constraints conflicting with overall volume/demand.`n`nDetailed Rule Evaluations:`nRule 1 ('# of Transitions'): The rule requires at least 1 transition for Bid ID 1. Note: Requiring at least one transition on an individual Bid ID is a very strict requirement; it forces a non-incumbent allocation even when data or economic factors might not support a transition. Consider applying this rule conditionally or relaxing the requirement to improve feasibility.`n`nPlease review supplier capacities, demand figures, and custom rule constraints for adjustments."
$notes.Range("A2").Value = $notesText

# ---------------------------------------------------------------------
# 3) LP Model sheet
# ---------------------------------------------------------------------
$lp = $wb.Worksheets.Item("LP Model")
$lpText = "\* Sourcing_with_MultiTier_Rebates_Discounts *\`nMinimize`nOBJ: S_A + S_B + S_C - rebate_A - rebate_B - rebate_C`nSubject To`nBaseSpend_A: S0_A - 0.5252 x_A_1 - 64 x_A_10 - 70 x_A_2 - 55 x_A_3 - 23 x_A_4`n - 54 x_A_5 - 42 x_A_6 - 23 x_A_7 - 75 x_A_8 - 97 x_A_9 = 0`nBaseSpend_B: S0_B - 10 x_B_1 - 13 x_B_10 - 70 x_B_2 - 65 x_B_3 - 75 x_B_4`n - 34 x_B_5 - 24 x_B_6 - 85 x_B_7 - 13 x_B_8 - 56 x_B_9 = 0`nBaseSpend_C: S0_C - 24 x_C_1 - 15 x_C_10 - 75 x_C_2 - 60 x_C_3 - 24 x_C_4`n - 44 x_C_6 - 42 x_C_7 - 24 x_C_8 - 32 x_C_9 = 0`nCapacity_B_Bid_ID_1: x_B_1 <= 100000000`nCapacity_B_Bid_ID_10: x_B_10 <= 100000000`nCapacity_B_Bid_ID_2: x_B_2 <= 100000000`nCapacity_B_Bid_ID_3: x_B_3 <= 100000000`nCapacity_B_Bid_ID_4: x_B_4 <= 100000000`nCapacity_B_Bid_ID_5: x_B_5 <= 100000000`nCapacity_B_Bid_ID_6: x_B_6 <= 100000000`nCapacity_B_Bid_ID_7: x_B_7 <= 100000000`nCapacity_B_Bid_ID_8: x_B_8 <= 100000000`nCapacity_B_Bid_ID_9: x_B_9 <= 100000000`nCapacity_C_Bid_ID_1: x_C_1 <= 100000000`nCapacity_C_Bid_ID_10: x_C_10 <= 100000000`nCapacity_C_Bid_ID_2: x_C_2 <= 100000000`nCapacity_C_Bid_ID_3: x_C_3 <= 100000000`nCapacity_C_Bid_ID_4: x_C_4 <= 100000000`nCapacity_C_Bid_ID_5: x_C_5 <= 100000000`nCapacity_C_Bid_ID_6: x_C_6 <= 100000000`nCapacity_C_Bid_ID_7: x_C_7 <= 100000000`nCapacity_C_Bid_ID_8: x_C_8 <= 100000000`nCapacity_C_Bid_ID_9: x_C_9 <= 100000000`nDemand_1: x_A_1 + x_B_1 + x_C_1 = 700`nDemand_10: x_A_10 + x_B_10 + x_C_10 = 13`nDemand_2: x_A_2 + x_B_2 + x_C_2 = 9000`nDemand_3: x_A_3 + x_B_3 + x_C_3 = 600`nDemand_4: x_A_4 + x_B_4 + x_C_4 = 5670`nDemand_5: x_A_5 + x_B_5 + x_C_5 = 45`nDemand_6: x_A_6 + x_B_6 + x_C_6 = 242`nDemand_7: x_A_7 + x_B_7 + x_C_7 = 664`nDemand_8: x_A_8 + x_B_8 + x_C_8 = 24`nDemand_9: x_A_9 + x_B_9 + x_C_9 = 232`nDiscountTierLower_A_0: d_A - 19400000000 z_discount_A_0 >= -19400000000`nDiscountTierLower_A_1: - 0.01 S0_A + d_A - 19400000000 z_discount_A_1`n >= -19400000000`nDiscountTierLower_B_0: d_B - 97000000000 z_discount_B_0 >= -97000000000`nDiscountTierLower_B_1: - 0.03 S0_B + d_B - 97000000000 z_discount_B_1`n >= -97000000000`nDiscountTierLower_C_0: d_C - 97000000000 z_discount_C_0 >= -97000000000`nDiscountTierLower_C_1: - 0.04 S0_C + d_C - 97000000000 z_discount_C_1`n >= -97000000000`nDiscountTierMax_A_0: 19400000000 z_discount_A_0 <= 19400001000`nDiscountTierMax_B_0: 97000000000 z_discount_B_0 <= 97000000500`nDiscountTierMax_C_0: 97000000000 z_discount_C_0 <= 97000000500`n_dummy: __dummy = 0`nDiscountTierMin_A_0: __dummy >= 0`nDiscountTierMin_A_1: x_A_1 + x_A_10 + x_A_3 + x_A_4 + x_A_8 + x_A_9`n - 1000 z_discount_A_1 >= 0`nDiscountTierMin_B_0: __dummy >= 0`nDiscountTierMin_B_1: x_B_2 + x_B_5 + x_B_6 + x_B_7 - 500 z_discount_B_1 >= 0`nDiscountTierMin_C_0: __dummy >= 0`nDiscountTierMin_C_1: x_C_1 + x_C_10 + x_C_3 + x_C_4 + x_C_8 + x_C_9`n - 500 z_discount_C_1 >= 0`nDiscountTierSelect_A: z_discount_A_0 + z_discount_A_1 = 1`nDiscountTierSelect_B: z_discount_B_0 + z_discount_B_1 = 1`nDiscountTierSelect_C: z_discount_C_0 + z_discount_C_1 = 1`nDiscountTierUpper_A_0: d_A + 19400000000 z_discount_A_0 <= 19400000000`nDiscountTierUpper_A_1: - 0.01 S0_A + d_A + 19400000000 z_discount_A_1`n <= 19400000000`nDiscountTierUpper_B_0: d_B + 97000000000 z_discount_B_0 <= 97000000000`nDiscountTierUpper_B_1: - 0.03 S0_B + d_B + 97000000000 z_discount_B_1`n <= 97000000000`nDiscountTierUpper_C_0: d_C + 97000000000 z_discount_C_0 <= 97000000000`nDiscountTierUpper_C_1: - 0.04 S0_C + d_C + 97000000000 z_discount_C_1`n <= 97000000000`nEffectiveSpend_A: - S0_A + S_A + d_A = 0`nEffectiveSpend_B: - S0_B + S_B + d_B = 0`nEffectiveSpend_C: - S0_C + S_C + d_C = 0`nNonBid_C_5: x_C_5 = 0`nRebateTierLower_A_0: rebate_A - 19400000000 y_rebate_A_0 >= -19400000000`nRebateTierLower_A_1: - 0.1 S_A + rebate_A - 19400000000 y_rebate_A_1`n >= -19400000000`nRebateTierLower_B_0: rebate_B - 97000000000 y_rebate_B_0 >= -97000000000`nRebateTierLower_B_1: - 0.05 S_B + rebate_B - 97000000000 y_rebate_B_1`n >= -97000000000`nRebateTierLower_C_0: rebate_C - 97000000000 y_rebate_C_0 >= -97000000000`nRebateTierLower_C_1: - 0.07 S_C + rebate_C - 97000000000 y_rebate_C_1`n >= -97000000000`nRebateTierMax_A_0: 19400000000 y_rebate_A_0 <= 19400000500`nRebateTierMax_B_0: 97000000000 y_rebate_B_0 <= 97000000500`nRebateTierMax_C_0: 97000000000 y_rebate_C_0 <= 97000000700`nRebateTierMin_A_0: __dummy >= 0`nRebateTierMin_A_1: - 500 y_rebate_A_1 >= 0`nRebateTierMin_B_0: __dummy >= 0`nRebateTierMin_B_1: x_B_2 + x_B_5 + x_B_6 + x_B_7 - 500 y_rebate_B_1 >= 0`nRebateTierMin_C_0: __dummy >= 0`nRebateTierMin_C_1: x_C_1 + x_C_10 + x_C_3 + x_C_4 + x_C_8 + x_C_9`n - 700 y_rebate_C_1 >= 0`nRebateTierSelect_A: y_rebate_A_0 + y_rebate_A_1 = 1`nRebateTierSelect_B: y_rebate_B_0 + y_rebate_B_1 = 1`nRebateTierSelect_C: y_rebate_C_0 + y_rebate_C_1 = 1`nRebateTierUpper_A_0: rebate_A + 19400000000 y_rebate_A_0 <= 19400000000`nRebateTierUpper_A_1: - 0.1 S_A + rebate_A + 19400000000 y_rebate_A_1`n <= 19400000000`nRebateTierUpper_B_0: rebate_B + 97000000000 y_rebate_B_0 <= 97000000000`nRebateTierUpper_B_1: - 0.05 S_B + rebate_B + 97000000000 y_rebate_B_1`n <= 97000000000`nRebateTierUpper_C_0: rebate_C + 97000000000 y_rebate_C_0 <= 97000000000`nRebateTierUpper_C_1: - 0.07 S_C + rebate_C + 97000000000 y_rebate_C_1`n <= 97000000000`nRule_0: __dummy >= 1`nTransition_10_A: - 13 T_10_A + x_A_10 <= 0`nTransition_10_B: - 13 T_10_B + x_B_10 <= 0`nTransition_1_B: - 700 T_1_B + x_B_1 <= 0`nTransition_1_C: - 700 T_1_C + x_C_1 <= 0`nTransition_2_A: - 9000 T_2_A + x_A_2 <= 0`nTransition_2_C: - 9000 T_2_C + x_C_2 <= 0`nTransition_3_A: - 600 T_3_A + x_A_3 <= 0`nTransition_3_B: - 600 T_3_B + x_B_3 <= 0`nTransition_4_A: - 5670 T_4_A + x_A_4 <= 0`nTransition_4_B: - 5670 T_4_B + x_B_4 <= 0`nTransition_5_A: - 45 T_5_A + x_A_5 <= 0`nTransition_5_B: - 45 T_5_B + x_B_5 <= 0`nTransition_6_A: - 242 T_6_A + x_A_6 <= 0`nTransition_6_B: - 242 T_6_B + x_B_6 <= 0`nTransition_7_A: - 664 T_7_A + x_A_7 <= 0`nTransition_7_B: - 664 T_7_B + x_B_7 <= 0`nTransition_8_A: - 24 T_8_A + x_A_8 <= 0`nTransition_8_B: - 24 T_8_B + x_B_8 <= 0`nTransition_9_A: - 232 T_9_A + x_A_9 <= 0`nTransition_9_B: - 232 T_9_B + x_B_9 <= 0`nVolume_A: V_A - x_A_1 - x_A_10 - x_A_2 - x_A_3 - x_A_4 - x_A_5 - x_A_6 - x_A_7`n - x_A_8 - x_A_9 = 0`nVolume_B: V_B - x_B_1 - x_B_10 - x_B_2 - x_B_3 - x_B_4 - x_B_5 - x_B_6 - x_B_7`n - x_B_8 - x_B_9 = 0`nVolume_C: V_C - x_C_1 - x_C_10 - x_C_2 - x_C_3 - x_C_4 - x_C_5 - x_C_6 - x_C_7`n - x_C_8 - x_C_9 = 0`nBinaries`nT_10_A`nT_10_B`nT_1_B`nT_1_C`nT_2_A`nT_2_C`nT_3_A`nT_3_B`nT_4_A`nT_4_B`nT_5_A`nT_5_B`nT_6_A`nT_6_B`nT_7_A`nT_7_B`nT_8_A`nT_8_B`nT_9_A`nT_9_B`ny_rebate_A_0`ny_rebate_A_1`ny_rebate_B_0`ny_rebate_B_1`ny_rebate_C_0`ny_rebate_C_1`nz_discount_A_0`nz_discount_A_1`nz_discount_B_0`nz_discount_B_1`nz_discount_C_0`nz_discount_C_1`nEnd`n"
$lp.Range("A2").Value = $lpText

Write-Output "Edit applied."
